# Update Name of Algo
# Applies refreshed imputation values produced by a re-run of the
# RandomForest algorithm (result_data_RandomForest.xlsx) on Sheet1.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = -7.700799999999999
$ws.Range("D5").Value = -8.054500000000008
$ws.Range("D6").Value = -8.1595
$ws.Range("B11").Value = 6.258999999999998
$ws.Range("A12").Value = -21.41779999999999
$ws.Range("B23").Value = 8.434800000000003
$ws.Range("C24").Value = -12.5385
$ws.Range("D27").Value = -7.918099999999999
$ws.Range("B28").Value = 5.663600000000002
$ws.Range("A32").Value = -21.1021
$ws.Range("B32").Value = 6.836199999999996
$ws.Range("B34").Value = 9.959800000000008
$ws.Range("A36").Value = -19.8952
$ws.Range("A38").Value = -19.80959999999999
$ws.Range("C38").Value = -10.8312
$ws.Range("B42").Value = 10.1315
$ws.Range("A46").Value = -21.97959999999999
$ws.Range("C52").Value = -11.035
$ws.Range("A54").Value = -21.90069999999999
$ws.Range("B54").Value = 4.7242
$ws.Range("A55").Value = -21.70970000000001
$ws.Range("D55").Value = -7.882800000000002
$ws.Range("A67").Value = -21.50739999999996
$ws.Range("A69").Value = -21.61879999999997
$ws.Range("A72").Value = -22.0325
$ws.Range("C78").Value = -13.2593
$ws.Range("D80").Value = -7.617499999999998
$ws.Range("C83").Value = -13.6642
$ws.Range("C85").Value = -14.26279999999999
$ws.Range("C86").Value = -14.14709999999999
$ws.Range("A91").Value = -20.89509999999997
$ws.Range("D95").Value = -7.438500000000004
$ws.Range("C96").Value = -10.3592
$ws.Range("B97").Value = 6.798199999999997
$ws.Range("D98").Value = -8.422600000000001
$ws.Range("A99").Value = -22.01109999999999
$ws.Range("B99").Value = 5.633499999999997
$ws.Range("B101").Value = 4.3532
$ws.Range("C103").Value = -13.62959999999999
$ws.Range("A104").Value = -21.7119
